$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two missing "Topic" entries for the class days that were
# added to the syllabus table (2024-02-17 / Saturday, row 7 and
# 2024-02-18 / Sunday, row 8).
$ws.Range("C7").Value = "Description list and text formatting in HTML"
$ws.Range("C8").Value = "Saturday: Hoiliday"

# The data rows grew a touch taller (18.75pt -> 19.5pt), and the last
# row grew along with them (19.5pt -> 20.25pt).
$ws.Range("A2:A13").EntireRow.RowHeight = 19.5
$ws.Range("A14").EntireRow.RowHeight = 20.25
